$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Agosto de 2020 a las 02:48"

# --- Country-name swaps in column A (use a temp placeholder per pair
#     to avoid transient shared-string collisions during the swap) ---

# Gabon <-> Libia (rows 94/95)
$ws.Range("A94").Value = "__TMP1__"
$ws.Range("A95").Value = "Gabon"
$ws.Range("A94").Value = "Libia"

# Letonia <-> Bahamas (rows 148/149)
$ws.Range("A148").Value = "__TMP2__"
$ws.Range("A149").Value = "Letonia"
$ws.Range("A148").Value = "Bahamas"

# Islas Caimanes <-> Polinesia Francesa (rows 185/186)
$ws.Range("A185").Value = "__TMP3__"
$ws.Range("A186").Value = "Islas Caimanes"
$ws.Range("A185").Value = "Polinesia Francesa"

# Rotation across rows 210-214:
# before: Groenlandia, Bonaire.., San Bartolome, Montserrat, Islas Malvinas
# after : San Bartolome, Groenlandia, Bonaire.., Islas Malvinas, Montserrat
$ws.Range("A210").Value = "__TMP4__"
$ws.Range("A211").Value = "__TMP5__"
$ws.Range("A212").Value = "__TMP6__"
$ws.Range("A213").Value = "__TMP7__"
$ws.Range("A214").Value = "__TMP8__"
$ws.Range("A210").Value = "San Bartolome"
$ws.Range("A211").Value = "Groenlandia"
$ws.Range("A212").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# --- Numeric data updates (columns B-H) ---
# row 4
$ws.Range("B4").Value = 5612011
$ws.Range("C4").Value = 40596
$ws.Range("D4").Value = 2970835
$ws.Range("E4").Value = 2467460
$ws.Range("G4").Value = 589
$ws.Range("H4").Value = 173716
# row 22
$ws.Range("B22").Value = 226686
$ws.Range("C22").Value = 1689
$ws.Range("E22").Value = 14490
# row 27
$ws.Range("B27").Value = 122872
$ws.Range("C27").Value = 785
$ws.Range("D27").Value = 109059
$ws.Range("E27").Value = 4781
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 9032
# row 81
$ws.Range("B81").Value = 14500
$ws.Range("C81").Value = 135
$ws.Range("D81").Value = 9442
$ws.Range("E81").Value = 4546
$ws.Range("G81").Value = 14
$ws.Range("H81").Value = 512
# row 94
$ws.Range("B94").Value = 8579
$ws.Range("C94").Value = 407
$ws.Range("D94").Value = 969
$ws.Range("E94").Value = 7453
$ws.Range("G94").Value = 4
$ws.Range("H94").Value = 157
# row 95
$ws.Range("B95").Value = 8270
$ws.Range("C95").Value = 45
$ws.Range("D95").Value = 6404
$ws.Range("E95").Value = 1813
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = 53
# row 148
$ws.Range("B148").Value = 1329
$ws.Range("C148").Value = 14
$ws.Range("D148").Value = 191
$ws.Range("E148").Value = 1119
$ws.Range("G148").Value = 1
$ws.Range("H148").Value = 19
# row 149
$ws.Range("B149").Value = 1323
$ws.Range("C149").Value = 1
$ws.Range("D149").Value = 1078
$ws.Range("E149").Value = 213
$ws.Range("H149").Value = 32
# row 160
$ws.Range("E160").Value = 645
$ws.Range("G160").Value = 5
$ws.Range("H160").Value = 30
# row 161
$ws.Range("D161").Value = 821
$ws.Range("E161").Value = 49
# row 171
$ws.Range("D171").Value = 336
$ws.Range("E171").Value = 76
# row 185
$ws.Range("B185").Value = 211
$ws.Range("C185").Value = 45
$ws.Range("D185").Value = 68
$ws.Range("E185").Value = 143
$ws.Range("H185").Value = 0
# row 186
$ws.Range("B186").Value = 203
$ws.Range("D186").Value = 202
$ws.Range("E186").Value = 0
$ws.Range("H186").Value = 1
# row 189
$ws.Range("B189").Value = 148
$ws.Range("C189").Value = 2
$ws.Range("E189").Value = 30
# row 193
$ws.Range("B193").Value = 121
$ws.Range("C193").Value = 12
$ws.Range("E193").Value = 67
$ws.Range("G193").Value = 1
$ws.Range("H193").Value = 5
# row 196
$ws.Range("B196").Value = 58
$ws.Range("C196").Value = 1
$ws.Range("E196").Value = 3
# row 210
$ws.Range("B210").Value = 16
$ws.Range("C210").Value = 3
$ws.Range("D210").Value = 9
$ws.Range("E210").Value = 7
# row 211
$ws.Range("B211").Value = 14
$ws.Range("D211").Value = 14
$ws.Range("E211").Value = 0
# row 212
$ws.Range("D212").Value = 7
$ws.Range("E212").Value = 6
# row 213
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0
# row 214
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
